$d = $word.ActiveDocument

# Locate the exact run of text that needs to be edited: "1.0"
# (the value of the "TC_IMAGE_LABEL" configuration key).
$rng = $d.Content
$found = $rng.Find.Execute('"1.0"', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the target text '""1.0""'"
}

$base = $rng.Start

# $rng now spans exactly the 5 characters:  "  1  .  0  "
#                                 offsets:  0  1  2  3  4
# The edit bumps the version from 1.0 to 1.1, i.e. the single
# character "0" (offset 3) is replaced by "1".
$zero = $d.Range($base + 3, $base + 4)
$zero.Text = "1"

# That in-place edit keeps everything as a single run (same
# formatting throughout), but the source document actually ends up
# with the run split into three sibling runs - "1. / 1 / " - all
# carrying identical run formatting. Reproduce that by nudging the
# formatting (apply then immediately revert Bold) on the two outer
# slices so the engine materializes them as their own runs instead of
# folding back into the neighbouring text.
$first = $d.Range($base, $base + 3)
$first.Font.Bold = 1
$first.Font.Bold = 0

$last = $d.Range($base + 4, $base + 5)
$last.Font.Bold = 1
$last.Font.Bold = 0

Write-Output "Result: [$($d.Range($base, $base + 5).Text)]"
